# MassWateR WQX Meta Template update
# - Remove "PON, " from the Parameter field instructions on the Instructions tab
# - Update "Available Values" for Result Sample Fraction with a recommended list
# - Change the "Field" entries for TP / TDP / E.coli on the Meta tab from
#   MassWateR to MassBays
# - Leave the Meta tab active/selected (it was "Instructions" before)

$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Meta")
$wsInstr = $wb.Worksheets.Item("Instructions")

# --- Instructions tab content edits -----------------------------------
# Parameter row instructions: drop "PON, " from the fraction example list
$wsInstr.Cells.Item(6, 2).Value = "Name of the measured parameter.`n- Note that this can be either the WQX or Simple parameter name.  However, if a parameter is distinguished by Sample Fraction only (i.e. TDP, TDN), then the Simple parameter name must be used here and in all other files (Results, DQO, etc.)"

# Result Sample Fraction row, Available Values column: add recommendation text
$wsInstr.Cells.Item(9, 4).Value = "standard list in WQX`nRecommended:  Filtered, lab; Filtered, field; Unfiltered; Non-Filterable (Particle)"
$wsInstr.Range("D9").WrapText = $true

# --- Meta tab content edits --------------------------------------------
# Sampling Method Context ("Field") changes from MassWateR to MassBays for
# TP, TDP, and E.coli rows
$wsMeta.Cells.Item(4, 2).Value = "MassBays"
$wsMeta.Cells.Item(5, 2).Value = "MassBays"
$wsMeta.Cells.Item(6, 2).Value = "MassBays"

# Left-align the Parameter column (A) on the Meta tab
$wsMeta.Range("A1:A6").HorizontalAlignment = -4131

# --- View / selection state ---------------------------------------------
# Make Meta the active/selected tab (it was Instructions before)
$wsInstr.Range("D10").Select()
$wsMeta.Range("B10").Select()
$wsMeta.Activate()
